$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I, row 2: blank cell with the same "thick-bottom-border" formatting
# as the rest of row 2 (H2).
$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial(-4122)

# Column I, row 3: new "2021" header, same formatting as the other year
# headers (H3) but bumped to 11pt (matches the new font/style added to the
# workbook).
$ws.Range("H3").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("I3").Value = 2021
$ws.Range("I3").Font.Size = 11

# Column I, row 4: new data value, same formatting as the rest of row 4
# (H4) but bumped to 11pt.
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = 149
$ws.Range("I4").Font.Size = 11

# Column I, row 5: new data value, same formatting as the rest of row 5
# (H5) but bumped to 11pt.
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").Value = 159
$ws.Range("I5").Font.Size = 11

# Match the author's final selection/cursor position recorded in the file.
[void]$ws.Range("K4").Select()
